$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows with changed data ---
# Note: D-column price values that are plain decimals (e.g. "373.17") are
# written with a leading apostrophe so Excel stores them as literal text
# (preserving trailing zeros etc.) instead of auto-converting to a number.
$ws.Range("D2").Value = "50.663.36"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "2.911.44"
$ws.Range("E3").Value = "  -2.24%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'373.17"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("D6").Value = "'99.20"
$ws.Range("E6").Value = "  -3.93%  "
$ws.Range("D7").Value = "'0.532"
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.579"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("D10").Value = "'35.61"
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "'0.0840"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "3.385.00"
$ws.Range("E13").Value = "  -1.88%  "
$ws.Range("D14").Value = "'17.82"
$ws.Range("D15").Value = "'7.49"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "2.910.72"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("D17").Value = "'10.99"
$ws.Range("E17").Value = "  +48.76%  "
$ws.Range("D18").Value = "'0.980"
$ws.Range("E18").Value = "  -2.37%  "
$ws.Range("D19").Value = "50.658.68"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").Value = "'3.02"
$ws.Range("E20").Value = "  -7.21%  "
$ws.Range("D21").Value = "'12.26"
$ws.Range("E21").Value = "  -4.39%  "
$ws.Range("D22").Value = "0.0₃0946"
$ws.Range("E22").Value = "  -1.24%  "
$ws.Range("D23").Value = "'68.63"
$ws.Range("E23").Value = "  -0.49%  "
$ws.Range("D24").Value = "'263.37"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").Value = "'3.11"
$ws.Range("E25").Value = "  +7.25%  "
$ws.Range("D26").Value = "'7.94"
$ws.Range("E26").Value = "  -2.64%  "
$ws.Range("D27").Value = "'7.30"
$ws.Range("E27").Value = "  -3.18%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'25.31"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").Value = "'0.161"
$ws.Range("E30").Value = "  -4.73%  "
$ws.Range("E31").Value = "  -8.38%  "
$ws.Range("D32").Value = "'9.90"
$ws.Range("D33").Value = "'50.74"
$ws.Range("E33").Value = "  -0.33%  "
$ws.Range("E34").Value = "  -0.91%  "
$ws.Range("D35").Value = "'32.73"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("D36").Value = "'0.0433"
$ws.Range("E36").Value = "  -4.66%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'3.05"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("D39").Value = "'0.114"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("D40").Value = "'16.29"
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("D41").Value = "'1.78"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("D44").Value = "'20.83"
$ws.Range("E44").Value = "  -3.05%  "
$ws.Range("D45").Value = "'2.05"
$ws.Range("E45").Value = "  -1.98%  "
$ws.Range("D46").Value = "'3.32"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "'0.265"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("D49").Value = "1.975.00"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("E50").Value = "  -2.95%  "
$ws.Range("D51").Value = "'5.15"
$ws.Range("E51").Value = "  +0.77%  "

# --- Rows 42/43: Monero and Stacks swap places (with updated price/volume) ---
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.42"
$ws.Range("E42").Value = "  -5.61%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "'118.88"
$ws.Range("E43").Value = "  -3.12%  "
